$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values that were repulled/recalculated
$ws.Range("F2").Value = -4
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = -8
$ws.Range("F6").Value = -4
$ws.Range("F10").Value = -1
